# Append newly scraped Lancers listings as of 2026-01-11 12:36:58.
# - Insert one new listing above the existing row (new row 2), pushing the
#   previously-existing listing down to row 3 (and refreshing its scrape
#   timestamp to the latest run).
# - Append one more newly scraped listing as a brand-new row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-01-11 12:36:58"

# Hyperlinks anchor themselves to absolute row/column positions and are not
# shifted automatically when rows are inserted, so drop them first and
# recreate them afterwards once all data is in its final place.
$ws.Hyperlinks.Delete()

# Shift the existing data row (row 2) down to row 3 to make room for the
# newest listing at the top.
$ws.Rows.Item(2).Insert()

# --- Row 2: new listing (inserted above the old top entry) ---
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【シンプル版】生成AIデジタル・コミュニティ制作の依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5469128"
$ws.Range("G2").Value = 310
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3: previously existing listing, timestamp refreshed ---
$ws.Range("A3").Value = $timestamp

# --- Row 4: new listing appended at the end ---
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "cordovaプロジェクトのバージョンアップ"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5469169"
$ws.Range("G4").Value = 10

# Recreate the hyperlinks (and their "Hyperlink" style) on the URL column
# for every data row now present.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5469128")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468866")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5469169")
$ws.Range("F4").Style = "Hyperlink"

# Widen column D (price) slightly to fit the longer new values.
$ws.Columns.Item(4).ColumnWidth = 28 - (5 / 6)
